# Auto-generated PowerShell Excel COM-interop script
# Applies the betexplorer scotland league-two 2023-2024 data refresh:
#  1) Re-shuffles odds/result data (cols F:V) among rows that share identical kickoff date
#     but had their betexplorer fixture order reshuffled by the scraping script.
#  2) Appends 7 new match rows (80-86) with the newly scraped fixtures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Step 1: permute columns F:V across rows that got reordered ----
$permutedRows = @{}
$permutedRows[12] = @("Spartans", 2, "Elgin City", 1, 1.73, "17/08/2023 09:13", 1.58, "19/08/2023 15:57", 3.62, "17/08/2023 09:13", 4.09, "19/08/2023 15:57", 4.03, "17/08/2023 09:13", 5.54, "19/08/2023 15:57", "https://www.betexplorer.com/football/scotland/league-two/spartans-elgin-city/dKsKy8Em/")
$permutedRows[13] = @("East Fife", 0, "Stenhousemuir", 2, 2.31, "17/08/2023 09:13", 2.32, "19/08/2023 15:55", 3.27, "17/08/2023 09:13", 3.52, "19/08/2023 15:55", 2.77, "17/08/2023 09:13", 2.93, "19/08/2023 15:55", "https://www.betexplorer.com/football/scotland/league-two/east-fife-stenhousemuir/vqmFxlas/")
$permutedRows[15] = @("Clyde", 1, "Peterhead", 2, 1.93, "17/08/2023 09:13", 3.12, "19/08/2023 14:24", 3.36, "17/08/2023 09:13", 3.65, "19/08/2023 14:24", 3.5, "17/08/2023 09:13", 2.16, "19/08/2023 14:24", "https://www.betexplorer.com/football/scotland/league-two/clyde-peterhead/rD1qrAyK/")
$permutedRows[16] = @("Bonnyrigg Rose", 1, "Stranraer", 1, 2.25, "17/08/2023 09:13", 2.25, "19/08/2023 15:58", 3.13, "17/08/2023 09:13", 3.23, "19/08/2023 15:58", 2.96, "17/08/2023 09:13", 3.29, "19/08/2023 15:58", "https://www.betexplorer.com/football/scotland/league-two/bonnyrigg-rose-stranraer/U9cuqjMD/")
$permutedRows[17] = @("Elgin City", 1, "East Fife", 1, 2.49, "24/08/2023 09:13", 2.44, "26/08/2023 15:59", 3.18, "24/08/2023 09:13", 3.5, "26/08/2023 15:58", 2.6, "24/08/2023 09:13", 2.76, "26/08/2023 15:59", "https://www.betexplorer.com/football/scotland/league-two/elgin-city-east-fife/zBtOzSTg/")
$permutedRows[18] = @("Forfar Athletic", 1, "Bonnyrigg Rose", 2, 2.12, "24/08/2023 09:13", 2.46, "26/08/2023 15:57", 3.17, "24/08/2023 09:13", 3.03, "26/08/2023 15:57", 3.17, "24/08/2023 09:13", 3.11, "26/08/2023 15:57", "https://www.betexplorer.com/football/scotland/league-two/forfar-athletic-bonnyrigg-rose/UXuSZoq0/")
$permutedRows[19] = @("Peterhead", 0, "Spartans", 1, 2.42, "24/08/2023 09:13", 2.6, "26/08/2023 15:55", 3.24, "24/08/2023 09:13", 3.42, "26/08/2023 15:55", 2.64, "24/08/2023 09:13", 2.63, "26/08/2023 15:55", "https://www.betexplorer.com/football/scotland/league-two/peterhead-spartans/SIbVY5b6/")
$permutedRows[20] = @("Stenhousemuir", 2, "Dumbarton", 4, 2.12, "24/08/2023 09:13", 2.69, "26/08/2023 15:54", 3.29, "24/08/2023 09:13", 3.19, "26/08/2023 15:54", 3.07, "24/08/2023 09:13", 2.69, "26/08/2023 15:54", "https://www.betexplorer.com/football/scotland/league-two/stenhousemuir-dumbarton/AslzYPDC/")
$permutedRows[21] = @("Stranraer", 1, "Clyde", 0, 1.65, "24/08/2023 09:13", 1.79, "26/08/2023 15:37", 3.75, "24/08/2023 09:13", 3.76, "26/08/2023 15:44", 4.52, "24/08/2023 09:13", 4.25, "26/08/2023 15:37", "https://www.betexplorer.com/football/scotland/league-two/stranraer-clyde/lYnvXqTI/")
$permutedRows[22] = @("Dumbarton", 1, "East Fife", 0, 1.9, "31/08/2023 09:13", 1.91, "02/09/2023 15:58", 3.41, "31/08/2023 09:13", 3.5, "02/09/2023 15:58", 3.52, "31/08/2023 09:13", 4.05, "02/09/2023 15:58", "https://www.betexplorer.com/football/scotland/league-two/dumbarton-east-fife/MLHhj7Tt/")
$permutedRows[23] = @("Clyde", 0, "Forfar Athletic", 0, 2.18, "31/08/2023 09:12", 3.23, "02/09/2023 15:59", 3.26, "31/08/2023 09:12", 3.48, "02/09/2023 15:59", 2.95, "31/08/2023 09:12", 2.17, "02/09/2023 15:59", "https://www.betexplorer.com/football/scotland/league-two/clyde-forfar-athletic/GvuQfmLP/")
$permutedRows[24] = @("Bonnyrigg Rose", 5, "Elgin City", 1, 1.75, "31/08/2023 09:13", 1.74, "02/09/2023 15:58", 3.52, "31/08/2023 09:13", 3.58, "02/09/2023 15:58", 4.08, "31/08/2023 09:13", 4.95, "02/09/2023 15:58", "https://www.betexplorer.com/football/scotland/league-two/bonnyrigg-rose-elgin-city/xltMeT5J/")
$permutedRows[25] = @("Spartans", 0, "Stenhousemuir", 1, 2.01, "31/08/2023 09:12", 2, "02/09/2023 15:59", 3.37, "31/08/2023 09:12", 3.59, "02/09/2023 15:59", 3.26, "31/08/2023 09:12", 3.57, "02/09/2023 15:59", "https://www.betexplorer.com/football/scotland/league-two/spartans-stenhousemuir/UoL0lobg/")
$permutedRows[37] = @("Dumbarton", 3, "Stranraer", 1, 1.85, "28/09/2023 08:13", 1.76, "30/09/2023 14:05", 3.43, "28/09/2023 08:13", 3.84, "30/09/2023 14:05", 3.71, "28/09/2023 08:13", 4.36, "30/09/2023 14:05", "https://www.betexplorer.com/football/scotland/league-two/dumbarton-stranraer/tr0u0rc5/")
$permutedRows[38] = @("East Fife", 1, "Forfar Athletic", 1, 2.18, "28/09/2023 08:13", 2.47, "30/09/2023 15:43", 3.15, "28/09/2023 08:13", 3.26, "30/09/2023 15:43", 3.08, "28/09/2023 08:13", 2.89, "30/09/2023 15:43", "https://www.betexplorer.com/football/scotland/league-two/east-fife-forfar-athletic/z5aqa2CB/")
$permutedRows[39] = @("Peterhead", 6, "Elgin City", 0, 1.58, "28/09/2023 08:13", 1.68, "30/09/2023 15:58", 3.85, "28/09/2023 08:13", 3.88, "30/09/2023 15:58", 4.98, "28/09/2023 08:13", 4.93, "30/09/2023 15:58", "https://www.betexplorer.com/football/scotland/league-two/peterhead-elgin-city/EebmbMRH/")
$permutedRows[40] = @("Spartans", 2, "Bonnyrigg Rose", 2, 1.97, "28/09/2023 08:13", 2.08, "30/09/2023 14:01", 3.33, "28/09/2023 08:13", 3.42, "30/09/2023 14:01", 3.4, "28/09/2023 08:13", 3.5, "30/09/2023 14:01", "https://www.betexplorer.com/football/scotland/league-two/spartans-bonnyrigg-rose/8dGLuEBj/")
$permutedRows[41] = @("Stenhousemuir", 2, "Clyde", 2, 1.63, "28/09/2023 08:13", 1.61, "30/09/2023 15:58", 3.72, "28/09/2023 08:13", 3.95, "30/09/2023 15:58", 4.54, "28/09/2023 08:13", 5.42, "30/09/2023 15:05", "https://www.betexplorer.com/football/scotland/league-two/stenhousemuir-clyde/U5HHtfdp/")
$permutedRows[42] = @("Forfar Athletic", 1, "Peterhead", 3, 2.44, "05/10/2023 08:12", 2.66, "07/10/2023 15:05", 3.02, "05/10/2023 08:12", 3.34, "07/10/2023 15:04", 2.77, "05/10/2023 08:12", 2.62, "07/10/2023 15:05", "https://www.betexplorer.com/football/scotland/league-two/forfar-athletic-peterhead/dIdG1eY9/")
$permutedRows[43] = @("Stranraer", 3, "Spartans", 4, 2.44, "05/10/2023 08:12", 2.67, "07/10/2023 15:27", 3.11, "05/10/2023 08:12", 3.41, "07/10/2023 15:04", 2.71, "05/10/2023 08:12", 2.57, "07/10/2023 15:27", "https://www.betexplorer.com/football/scotland/league-two/stranraer-spartans/jgBJ0FmG/")
$permutedRows[44] = @("Clyde", 0, "Dumbarton", 4, 3.4, "05/10/2023 08:12", 3.66, "07/10/2023 15:04", 3.3, "05/10/2023 08:12", 3.61, "07/10/2023 15:04", 1.98, "05/10/2023 08:12", 1.97, "07/10/2023 15:04", "https://www.betexplorer.com/football/scotland/league-two/clyde-dumbarton/z3873H3c/")
$permutedRows[54] = @("Stenhousemuir", 3, "Spartans", 2, 2.12, "09/11/2023 09:12", 2.37, "11/11/2023 15:58", 3.33, "09/11/2023 09:12", 3.61, "11/11/2023 15:58", 3.04, "09/11/2023 09:12", 2.79, "11/11/2023 15:58", "https://www.betexplorer.com/football/scotland/league-two/stenhousemuir-spartans/WC8p5Uv2/")
$permutedRows[58] = @("Forfar Athletic", 1, "Clyde", 1, 1.73, "09/11/2023 09:12", 1.89, "11/11/2023 15:49", 3.59, "09/11/2023 09:12", 3.64, "11/11/2023 15:49", 4.08, "09/11/2023 09:12", 3.92, "11/11/2023 15:55", "https://www.betexplorer.com/football/scotland/league-two/forfar-athletic-clyde/zoFy7jAk/")
$permutedRows[60] = @("Stranraer", 0, "Stenhousemuir", 3, 2.95, "16/11/2023 09:13", 3.78, "18/11/2023 15:18", 3.33, "16/11/2023 09:13", 3.49, "18/11/2023 15:58", 2.17, "16/11/2023 09:13", 1.97, "18/11/2023 15:18", "https://www.betexplorer.com/football/scotland/league-two/stranraer-stenhousemuir/r9BXQB1L/")
$permutedRows[61] = @("Bonnyrigg Rose", 0, "Forfar Athletic", 2, 2.03, "16/11/2023 09:13", 2.58, "18/11/2023 15:58", 3.33, "16/11/2023 09:13", 3.09, "18/11/2023 15:52", 3.23, "16/11/2023 09:13", 2.88, "18/11/2023 15:58", "https://www.betexplorer.com/football/scotland/league-two/bonnyrigg-rose-forfar-athletic/pGdSTDH2/")
$permutedRows[62] = @("Clyde", 0, "East Fife", 4, 2.85, "16/11/2023 09:13", 2.86, "18/11/2023 15:32", 3.24, "16/11/2023 09:13", 3.25, "18/11/2023 15:32", 2.27, "16/11/2023 09:13", 2.5, "18/11/2023 15:32", "https://www.betexplorer.com/football/scotland/league-two/clyde-east-fife/O6eWSXW8/")
$permutedRows[64] = @("Dumbarton", 1, "Elgin City", 0, 1.44, "16/11/2023 09:13", 1.38, "18/11/2023 15:52", 4.29, "16/11/2023 09:13", 4.86, "18/11/2023 15:52", 5.68, "16/11/2023 09:13", 7.89, "18/11/2023 15:52", "https://www.betexplorer.com/football/scotland/league-two/dumbarton-elgin-city/4ICTRinF/")
$permutedRows[65] = @("Bonnyrigg Rose", 4, "East Fife", 2, 2.05, "05/10/2023 09:12", 2.65, "21/11/2023 17:48", 3.26, "05/10/2023 09:12", 3.07, "21/11/2023 18:47", 3.26, "05/10/2023 09:12", 2.47, "21/11/2023 14:40", "https://www.betexplorer.com/football/scotland/league-two/bonnyrigg-rose-east-fife/Ma434cli/")
$permutedRows[66] = @("Elgin City", 1, "Forfar Athletic", 0, 2.79, "19/10/2023 09:13", 2.99, "21/11/2023 20:39", 3.17, "19/10/2023 09:13", 3.19, "21/11/2023 20:36", 2.34, "19/10/2023 09:13", 2.44, "21/11/2023 20:39", "https://www.betexplorer.com/football/scotland/league-two/elgin-city-forfar-athletic/lUjffDQq/")
$permutedRows[74] = @("Elgin City", 0, "Spartans", 4, 3.66, "14/12/2023 09:13", 4.68, "16/12/2023 15:53", 3.42, "14/12/2023 09:13", 3.71, "16/12/2023 15:57", 1.89, "14/12/2023 09:13", 1.74, "16/12/2023 15:52", "https://www.betexplorer.com/football/scotland/league-two/elgin-city-spartans/zLVpYS8l/")
$permutedRows[75] = @("Forfar Athletic", 2, "Dumbarton", 4, 3.13, "14/12/2023 09:13", 3.01, "16/12/2023 15:57", 3.17, "14/12/2023 09:13", 2.98, "16/12/2023 15:57", 2.19, "14/12/2023 09:13", 2.57, "16/12/2023 15:57", "https://www.betexplorer.com/football/scotland/league-two/forfar-athletic-dumbarton/UsVlXnOf/")
$permutedRows[76] = @("Peterhead", 2, "Clyde", 1, 1.55, "14/12/2023 09:13", 1.54, "16/12/2023 15:31", 4.01, "14/12/2023 09:13", 4.29, "16/12/2023 15:31", 5.12, "14/12/2023 09:13", 5.75, "16/12/2023 15:31", "https://www.betexplorer.com/football/scotland/league-two/peterhead-clyde/roZhW6w1/")
$permutedRows[77] = @("Stenhousemuir", 2, "East Fife", 1, 1.78, "14/12/2023 09:13", 1.92, "16/12/2023 15:56", 3.53, "14/12/2023 09:13", 3.46, "16/12/2023 15:58", 4.01, "14/12/2023 09:13", 4.07, "16/12/2023 15:58", "https://www.betexplorer.com/football/scotland/league-two/stenhousemuir-east-fife/Q9OcVQh7/")
$permutedRows[78] = @("Stranraer", 3, "Bonnyrigg Rose", 1, 2.34, "14/12/2023 09:13", 2.67, "16/12/2023 15:59", 3.16, "14/12/2023 09:13", 3.26, "16/12/2023 15:59", 2.86, "14/12/2023 09:13", 2.66, "16/12/2023 15:59", "https://www.betexplorer.com/football/scotland/league-two/stranraer-bonnyrigg-rose/40N1Up8D/")

foreach ($rowNum in $permutedRows.Keys) {
    $rowVals = $permutedRows[$rowNum]
    $arr = New-Object 'object[,]' 1,17
    for ($i = 0; $i -lt 17; $i++) { $arr[0,$i] = $rowVals[$i] }
    $ws.Range("F" + $rowNum + ":V" + $rowNum).Value = $arr
}

# ---- Step 2: append the new rows (80-86), copying formats from row 79 ----
$ws.Range("A79:V79").Copy()
$ws.Range("A80:V86").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$newRows = @{}
$newRows[80] = @(79, 45283.66666666666, "Peterhead", 1, "Forfar Athletic", 2, 1.85, "21/12/2023 09:12", 2.02, "23/12/2023 15:57", 3.36, "21/12/2023 09:12", 3.28, "23/12/2023 15:57", 3.91, "21/12/2023 09:12", 3.87, "23/12/2023 15:57", "https://www.betexplorer.com/football/scotland/league-two/peterhead-forfar-athletic/jeOHzMUn/")
$newRows[81] = @(80, 45283.66666666666, "Dumbarton", 4, "Clyde", 4, 1.48, "21/12/2023 09:12", 1.53, "23/12/2023 15:53", 4.2, "21/12/2023 09:12", 4.22, "23/12/2023 15:53", 5.65, "21/12/2023 09:12", 6.01, "23/12/2023 15:53", "https://www.betexplorer.com/football/scotland/league-two/dumbarton-clyde/W6BqgUpE/")
$newRows[82] = @(81, 45283.66666666666, "Spartans", 3, "Stranraer", 0, 1.75, "21/12/2023 09:12", 1.72, "23/12/2023 08:12", 3.58, "21/12/2023 09:12", 3.83, "23/12/2023 14:00", 4.17, "21/12/2023 09:12", 4.61, "23/12/2023 14:00", "https://www.betexplorer.com/football/scotland/league-two/spartans-stranraer/bLRPYaaa/")
$newRows[83] = @(82, 45283.66666666666, "Stenhousemuir", 2, "Elgin City", 0, 1.36, "21/12/2023 09:12", 1.33, "23/12/2023 15:55", 4.54, "21/12/2023 09:12", 4.85, "23/12/2023 15:56", 7.32, "21/12/2023 09:12", 10.25, "23/12/2023 15:56", "https://www.betexplorer.com/football/scotland/league-two/stenhousemuir-elgin-city/4bSLZupg/")
$newRows[84] = @(83, 45290.66666666666, "Elgin City", 2, "Peterhead", 1, 4.5, "28/12/2023 09:12", 4.69, "30/12/2023 15:47", 3.59, "28/12/2023 09:12", 3.9, "30/12/2023 15:47", 1.69, "28/12/2023 09:12", 1.7, "30/12/2023 15:47", "https://www.betexplorer.com/football/scotland/league-two/elgin-city-peterhead/YwKxWcqI/")
$newRows[85] = @(84, 45290.66666666666, "Clyde", 1, "Stenhousemuir", 2, 4.41, "28/12/2023 09:12", 4.45, "30/12/2023 15:52", 3.64, "28/12/2023 09:12", 4, "30/12/2023 15:52", 1.69, "28/12/2023 09:12", 1.71, "30/12/2023 15:52", "https://www.betexplorer.com/football/scotland/league-two/clyde-stenhousemuir/fmLYWwUB/")
$newRows[86] = @(85, 45290.66666666666, "Forfar Athletic", 0, "East Fife", 0, 2.27, "28/12/2023 09:12", 2.29, "30/12/2023 15:53", 3.14, "28/12/2023 09:12", 3.28, "30/12/2023 15:53", 2.99, "28/12/2023 09:12", 3.17, "30/12/2023 15:53", "https://www.betexplorer.com/football/scotland/league-two/forfar-athletic-east-fife/S0JtVHbO/")

$newRowCols = @('A','E','F','G','H','I','J','K','L','M','N','O','P','Q','R','S','T','U','V')
foreach ($rowNum in $newRows.Keys) {
    $rowVals = $newRows[$rowNum]
    $ws.Cells.Item($rowNum, 2).Value = "scotland"
    $ws.Cells.Item($rowNum, 3).Value = "league-two"
    $ws.Cells.Item($rowNum, 4).Value = "2023-2024"
    for ($i = 0; $i -lt $newRowCols.Count; $i++) {
        $ws.Range($newRowCols[$i] + $rowNum).Value = $rowVals[$i]
    }
}

Write-Host "Applied betexplorer data refresh."
